$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 140.71428
$ws.Range("I9").Value = 131.6
$ws.Range("K9").Value = 131.6
$ws.Range("M9").Value = 37.40000000000001
$ws.Range("H17").Value = 721.5
$ws.Range("J17").Value = 721.5
$ws.Range("L17").Value = 2164.5
$ws.Range("N17").Value = -2500.5
$ws.Range("H38").Value = 817.1
$ws.Range("I38").Value = 146.375
$ws.Range("J38").Value = 3500
$ws.Range("K38").Value = 439.125
$ws.Range("L38").Value = 10500
$ws.Range("M38").Value = -67.125
$ws.Range("N38").Value = -11244
$ws.Range("H41").Value = 433.33334
$ws.Range("I41").Value = 480
$ws.Range("K41").Value = 480
$ws.Range("M41").Value = -40
$ws.Range("H58").Value = 4357.6665
$ws.Range("I58").Value = 1250
$ws.Range("J58").Value = 5245.5713
$ws.Range("K58").Value = 3750
$ws.Range("L58").Value = 15736.7139
$ws.Range("M58").Value = -3600
$ws.Range("N58").Value = -16036.7139
$ws.Range("H87").Value = 75199.7
$ws.Range("J87").Value = 75199.7
$ws.Range("L87").Value = 75199.7
$ws.Range("N87").Value = -77695.7
$ws.Range("H90").Value = 75199.7
$ws.Range("J90").Value = 75199.7
$ws.Range("L90").Value = 225599.1
$ws.Range("N90").Value = -238079.1
$ws.Range("H96").Value = 166669250
$ws.Range("I96").Value = 166669250
$ws.Range("K96").Value = 500007750
$ws.Range("M96").Value = -500006377
$ws.Range("H132").Value = 2622.2856
$ws.Range("I132").Value = 1177.44
$ws.Range("K132").Value = 3532.32
$ws.Range("M132").Value = -1002.32
$ws.Range("H138").Value = 2352.3967
$ws.Range("I138").Value = 1751.2106
$ws.Range("J138").Value = 2612
$ws.Range("K138").Value = 5253.6318
$ws.Range("L138").Value = 7836
$ws.Range("M138").Value = -113.6318000000001
$ws.Range("N138").Value = -18116

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 14493968
$ws.Range("I2").Value = 18519206
$ws.Range("J2").Value = 3109.2
$ws.Range("K2").Value = 18519206
$ws.Range("L2").Value = 3109.2
$ws.Range("M2").Value = -18519093
$ws.Range("N2").Value = -3335.2
$ws.Range("H28").Value = 2794
$ws.Range("I28").Value = 2794
$ws.Range("K28").Value = 2794
$ws.Range("M28").Value = -2602
$ws.Range("H99").Value = 2794
$ws.Range("I99").Value = 2794
$ws.Range("K99").Value = 2794
$ws.Range("M99").Value = 201
$ws.Range("H110").Value = 5293650.5
$ws.Range("I110").Value = 12347285
$ws.Range("J110").Value = 3425
$ws.Range("K110").Value = 12347285
$ws.Range("L110").Value = 3425
$ws.Range("M110").Value = -12345240
$ws.Range("N110").Value = -7515
$ws.Range("H116").Value = 14493968
$ws.Range("I116").Value = 18519206
$ws.Range("J116").Value = 3109.2
$ws.Range("K116").Value = 18519206
$ws.Range("L116").Value = 3109.2
$ws.Range("M116").Value = -18516912
$ws.Range("N116").Value = -7697.2

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 14493968
$ws.Range("I3").Value = 18519206
$ws.Range("J3").Value = 3109.2
$ws.Range("K3").Value = 18519206
$ws.Range("L3").Value = 3109.2
$ws.Range("M3").Value = -18519092
$ws.Range("N3").Value = -3337.2
$ws.Range("H105").Value = 4390343.5
$ws.Range("I105").Value = 5955980.5
$ws.Range("K105").Value = 5955980.5
$ws.Range("M105").Value = -5954233.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()  # was -82258
$ws.Range("H122").Value = 1431.8889
$ws.Range("I122").Value = 1426.7142
$ws.Range("K122").Value = 4280.142599999999
$ws.Range("M122").Value = -1830.142599999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 12000
$ws.Range("J3").Value = 12000
$ws.Range("L3").Value = 36000
$ws.Range("N3").Value = -36224
$ws.Range("H4").Value = 34393616
$ws.Range("I4").Value = 42036444
$ws.Range("J4").Value = 873.1667
$ws.Range("K4").Value = 126109332
$ws.Range("L4").Value = 2619.5001
$ws.Range("M4").Value = -126109220
$ws.Range("N4").Value = -2843.5001
$ws.Range("H86").Value = 250.66667
$ws.Range("I86").Value = 249
$ws.Range("K86").Value = 747
$ws.Range("M86").Value = 439
$ws.Range("H89").Value = 250.66667
$ws.Range("I89").Value = 249
$ws.Range("K89").Value = 2241
$ws.Range("M89").Value = 3687
$ws.Range("H140").Value = 1572.3334
$ws.Range("I140").Value = 1260.909
$ws.Range("J140").Value = 4998
$ws.Range("K140").Value = 3782.727
$ws.Range("L140").Value = 14994
$ws.Range("M140").Value = 1397.273
$ws.Range("N140").Value = -25354

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3496.1428
$ws.Range("I80").Value = 2987
$ws.Range("J80").Value = 3699.8
$ws.Range("K80").Value = 2987
$ws.Range("L80").Value = 3699.8
$ws.Range("M80").Value = -1989
$ws.Range("N80").Value = -5695.8
$ws.Range("H83").Value = 3496.1428
$ws.Range("I83").Value = 2987
$ws.Range("J83").Value = 3699.8
$ws.Range("K83").Value = 14935
$ws.Range("L83").Value = 18499
$ws.Range("M83").Value = -9943
$ws.Range("N83").Value = -28483
$ws.Range("H97").Value = 1540.7273
$ws.Range("J97").Value = 1815.6666
$ws.Range("L97").Value = 1815.6666
$ws.Range("N97").Value = -2807.6666
$ws.Range("H113").Value = 15919.6
$ws.Range("I113").Value = 4866
$ws.Range("J113").Value = 32500
$ws.Range("K113").Value = 4866
$ws.Range("L113").Value = 32500
$ws.Range("M113").Value = -2696
$ws.Range("N113").Value = -36840
$ws.Range("H122").Value = 204993.4
$ws.Range("I122").Value = 3999.5
$ws.Range("K122").Value = 11998.5
$ws.Range("M122").Value = -9548.5
$ws.Range("H126").Value = 14000
$ws.Range("J126").Value = 14000
$ws.Range("L126").Value = 42000
$ws.Range("N126").Value = -46940

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4973.375
$ws.Range("I40").Value = 4777.4
$ws.Range("K40").Value = 4777.4
$ws.Range("M40").Value = -4641.4
$ws.Range("H61").Value = 13892026
$ws.Range("I61").Value = 15876244
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 15876244
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -15876042
$ws.Range("N61").Value = -2904
$ws.Range("H93").Value = 1949
$ws.Range("H113").Value = 13892026
$ws.Range("I113").Value = 15876244
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 15876244
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = -15874074
$ws.Range("N113").Value = -6840

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1708
$ws.Range("I96").Value = 1708
$ws.Range("K96").Value = 1708
$ws.Range("M96").Value = -335
$ws.Range("H113").Value = 1506.75
$ws.Range("J113").Value = 2183.4285
$ws.Range("L113").Value = 6550.2855
$ws.Range("N113").Value = -10890.2855
